$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-17 (Sending cluster x Target cluster, full 4x4
# combination of ECs/FAPs/M2/sCs), columns A-T.
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
#          E=Ligand-expressing cells, F=Ligand detection rate,
#          G=Ligand average expression value, H=Ligand total expression value,
#          I=Ligand derived specificity (avg), J=Ligand derived specificity (total),
#          K=Receptor-expressing cells, L=Receptor detection rate,
#          M=Receptor average expression value, N=Receptor total expression value,
#          O=Receptor derived specificity (avg), P=Receptor derived specificity (total),
#          Q=Edge average expression weight, R=Edge total expression weight,
#          S=Edge average expression derived specificity, T=Edge total expression derived specificity
$csv = @"
ECs,Col1a1,Itga11,ECs,3,1,9.108069666666667,27.324209,0.00155006418458712,0.00155006418458712,1,0.3333333333333333,0.09233266666666667,0.276998,0.002874858548413657,0.002874858548413658,0.8409723605091112,7.568751244582001,4.456215271650126E-06,4.456215271650126E-06
ECs,Col1a1,Itga11,FAPs,3,1,9.108069666666667,27.324209,0.00155006418458712,0.00155006418458712,3,1,31.51785566666667,94.553567,0.9813360759751099,0.98133607597511,287.0668251559448,2583.601426403503,0.001521133904412282,0.001521133904412283
ECs,Col1a1,Itga11,M2,3,1,9.108069666666667,27.324209,0.00155006418458712,0.00155006418458712,1,0.3333333333333333,0.02446166666666667,0.07338500000000001,0.0007616354434881705,0.0007616354434881706,0.2227985641627778,2.005187077465,1.18058382266314E-06,1.18058382266314E-06
ECs,Col1a1,Itga11,sCs,3,1,9.108069666666667,27.324209,0.00155006418458712,0.00155006418458712,3,1,0.4826403333333333,1.447921,0.01502743003298815,0.01502743003298815,4.395921779943222,39.563296019489,2.329348108052376E-05,2.329348108052377E-05
FAPs,Col1a1,Itga11,ECs,3,1,5771.873535333333,17315.620606,0.9822909543423312,0.9822909543423313,1,0.3333333333333333,0.09233266666666667,0.276998,0.002874858548413657,0.002874858548413658,532.9324751800876,4796.392276620789,0.00282394754712046,0.002823947547120461
FAPs,Col1a1,Itga11,FAPs,3,1,5771.873535333333,17315.620606,0.9822909543423312,0.9822909543423313,3,1,31.51785566666667,94.553567,0.9813360759751099,0.98133607597511,181917.0770128891,1637253.693116002,0.9639575506001491,0.9639575506001493
FAPs,Col1a1,Itga11,M2,3,1,5771.873535333333,17315.620606,0.9822909543423312,0.9822909543423313,1,0.3333333333333333,0.02446166666666667,0.07338500000000001,0.0007616354434881705,0.0007616354434881706,141.1896464634789,1270.70681817131,0.0007481476066449397,0.0007481476066449399
FAPs,Col1a1,Itga11,sCs,3,1,5771.873535333333,17315.620606,0.9822909543423312,0.9822909543423313,3,1,0.4826403333333333,1.447921,0.01502743003298815,0.01502743003298815,2785.738967051125,25071.65070346013,0.01476130858841654,0.01476130858841654
M2,Col1a1,Itga11,ECs,3,1,1.272029666666667,3.816089,0.0002164813950916887,0.0002164813950916887,1,0.3333333333333333,0.09233266666666667,0.276998,0.002874858548413657,0.002874858548413658,0.1174498912024445,1.057049020822,6.223533892518557E-07,6.223533892518558E-07
M2,Col1a1,Itga11,FAPs,3,1,1.272029666666667,3.816089,0.0002164813950916887,0.0002164813950916887,3,1,31.51785566666667,94.553567,0.9813360759751099,0.98133607597511,40.09164743771811,360.824826939463,0.0002124410027808952,0.0002124410027808952
M2,Col1a1,Itga11,M2,3,1,1.272029666666667,3.816089,0.0002164813950916887,0.0002164813950916887,1,0.3333333333333333,0.02446166666666667,0.07338500000000001,0.0007616354434881705,0.0007616354434881706,0.03111596569611111,0.280043691265,1.648799033575962E-07,1.648799033575962E-07
M2,Col1a1,Itga11,sCs,3,1,1.272029666666667,3.816089,0.0002164813950916887,0.0002164813950916887,3,1,0.4826403333333333,1.447921,0.01502743003298815,0.01502743003298815,0.6139328223298889,5.525395400969,3.253159018184016E-06,3.253159018184016E-06
sCs,Col1a1,Itga11,ECs,3,1,93.67702500000001,281.031075,0.01594250007799006,0.01594250007799006,1,0.3333333333333333,0.09233266666666667,0.276998,0.002874858548413657,0.002874858548413658,8.649449523650002,77.84504571285002,4.583243263229513E-05,4.583243263229513E-05
sCs,Col1a1,Itga11,FAPs,3,1,93.67702500000001,281.031075,0.01594250007799006,0.01594250007799006,3,1,31.51785566666667,94.553567,0.9813360759751099,0.98133607597511,2952.498953232725,26572.49057909453,0.01564495046776765,0.01564495046776765
sCs,Col1a1,Itga11,M2,3,1,93.67702500000001,281.031075,0.01594250007799006,0.01594250007799006,1,0.3333333333333333,0.02446166666666667,0.07338500000000001,0.0007616354434881705,0.0007616354434881706,2.291496159875001,20.62346543887501,1.214237311721015E-05,1.214237311721015E-05
sCs,Col1a1,Itga11,sCs,3,1,93.67702500000001,281.031075,0.01594250007799006,0.01594250007799006,3,1,0.4826403333333333,1.447921,0.01502743003298815,0.01502743003298815,45.212310571675,406.910795145075,0.0002395748044729037,0.0002395748044729038
"@

$lines = $csv -split "`n"
$row = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $f = $line -split ","

    $ws.Cells.Item($row, 1).Value = $f[0]
    $ws.Cells.Item($row, 2).Value = $f[1]
    $ws.Cells.Item($row, 3).Value = $f[2]
    $ws.Cells.Item($row, 4).Value = $f[3]

    for ($col = 5; $col -le 20; $col++) {
        $ws.Cells.Item($row, $col).Value = [double]$f[$col - 1]
    }

    $row++
}
